# Add files via upload - added population data for ethnicity
$wb = $excel.ActiveWorkbook

# 1. ethnicities sheet: re-enter the B7:G7 delta formula as one pass so it
#    collapses into a shared formula group (matches the authoring diff).
$wsEth = $wb.Worksheets.Item("ethnicities")
$wsEth.Range("B7:G7").Formula = "=B5-B6"

# 2. prop sheet: append the new "Population" row (row 8) with the
#    ethnicity population percentages.
$wsProp = $wb.Worksheets.Item("prop")
$wsProp.Range("A8").Value = "Population"
$wsProp.Range("B8").Value = 76.8
$wsProp.Range("C8").Value = 16.2
$wsProp.Range("D8").Value = 9.2
$wsProp.Range("E8").Value = 3.7
$wsProp.Range("F8").Value = 3.2
$wsProp.Range("G8").Value = 0

# Column A on "prop" widens to fit the new label; re-fit it.
$wsProp.Columns.Item(1).AutoFit() | Out-Null

# Leave the cursor where the author's session ended up (just past the new row).
$wsProp.Range("G9").Select() | Out-Null
